$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 304
$ws.Range("F3").Value = 1068
$ws.Range("F5").Value = 1107
$ws.Range("F14").Value = 126
$ws.Range("F19").Value = 309
$ws.Range("F20").Value = 12
$ws.Range("F23").Value = 359
$ws.Range("F25").Value = 606
$ws.Range("F26").Value = 75952
$ws.Range("F27").Value = 75953
$ws.Range("F30").Value = 33068
$ws.Range("F31").Value = 33068
$ws.Range("F32").Value = 452
$ws.Range("F33").Value = 10
$ws.Range("F36").Value = 6
$ws.Range("F43").Value = 723
$ws.Range("F44").Value = 427
$ws.Range("F47").Value = 334
$ws.Range("F49").Value = 3
# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F15").Value = 941
$ws.Range("F18").Value = 395
$ws.Range("F35").Value = 1346
$ws.Range("F38").Value = 95
$ws.Range("F39").Value = 95
# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 552
# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 304
$ws.Range("F5").Value = 1068
$ws.Range("F8").Value = 1107
$ws.Range("F14").Value = 552
$ws.Range("F15").Value = 552
$ws.Range("F23").Value = 309
$ws.Range("F28").Value = 395
$ws.Range("F29").Value = 359
$ws.Range("F30").Value = 606
$ws.Range("F33").Value = 75953
$ws.Range("F35").Value = 33068
$ws.Range("F36").Value = 10
$ws.Range("F46").Value = 723
$ws.Range("F48").Value = 427
$ws.Range("F49").Value = 95
$ws.Range("F51").Value = 334

